# "thay doi nhan excel nguoi lao dong"
# Replace the note used on every Excel-date field ("ngaysinh", "bdhopdong",
# "kthopdong", "bddochai", "ktdochai", "bdbhxh", "ktbhxh") in column C.
# The old note told users to rely on Excel's native Date format; the new
# note instead tells them to enter the date as plain text in dd/mm/yyyy
# format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newNote = "Định dạng text. Ghi theo mẫu: dd/mm/yyyy (vd:01/01/2023)"

$ws.Range("C4").Value  = $newNote   # ngaysinh
$ws.Range("C18").Value = $newNote   # bdhopdong
$ws.Range("C19").Value = $newNote   # kthopdong
$ws.Range("C26").Value = $newNote   # bddochai
$ws.Range("C27").Value = $newNote   # ktdochai
$ws.Range("C30").Value = $newNote   # bdbhxh
$ws.Range("C31").Value = $newNote   # ktbhxh

# Match the saved cursor/selection position from the authored workbook.
$ws.Range("C31").Select()
